# Austria Bundesliga workbook update (17-02-2024 refresh).
# The source data rows were re-keyed upstream; several rows simply swapped
# places with a neighbouring row (everything except the leading "id" column
# A moved with the match), one trio rotated, and a handful of still-unplayed
# fixture rows had their opening-line odds refreshed in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# Simple two-row swaps (row "A" id stays put, all other columns trade places)
$swapPairs = @(
    @(17, 18),
    @(30, 31),
    @(89, 90),
    @(93, 94),
    @(100, 101),
    @(107, 108),
    @(114, 115),
    @(143, 144),
    @(159, 160),
    @(170, 171),
    @(203, 204)
)

foreach ($pair in $swapPairs) {
    Swap-Rows $pair[0] $pair[1]
}

# Three-way rotation: row 152 takes row 154's data, row 153 takes row 152's
# (original) data, row 154 takes row 153's (original) data.
$v152 = $ws.Range("B152:AC152").Value()
$v153 = $ws.Range("B153:AC153").Value()
$v154 = $ws.Range("B154:AC154").Value()

$ws.Range("B152:AC152").Value = $v154
$ws.Range("B153:AC153").Value = $v152
$ws.Range("B154:AC154").Value = $v153

# In-place odds refresh for not-yet-played fixtures (no row reshuffle here).
$ws.Range("N206").Value = 1.8
$ws.Range("O206").Value = 3.6
$ws.Range("P206").Value = 4.5
$ws.Range("Q206").Value = -0.5
$ws.Range("R206").Value = 1.85
$ws.Range("S206").Value = 2
$ws.Range("U206").Value = 2.025
$ws.Range("V206").Value = 1.825

$ws.Range("R207").Value = 1.925
$ws.Range("S207").Value = 1.925
$ws.Range("U207").Value = 1.925
$ws.Range("V207").Value = 1.925

$ws.Range("N208").Value = 8.5
$ws.Range("O208").Value = 4
$ws.Range("R208").Value = 1.775
$ws.Range("S208").Value = 2.1
$ws.Range("U208").Value = 2
$ws.Range("V208").Value = 1.85

$ws.Range("U209").Value = 1.85
$ws.Range("V209").Value = 2

$ws.Range("R210").Value = 1.95
$ws.Range("S210").Value = 1.9
